# "add webtest load validation"
#
# Adds a "Run" flag column (values "Y") to the two sheets that drive the
# webtest runner: TestSteps (new column F) and PageModules (new column I).
# Formats are picked up from the sheet's existing header/data cell to the
# left so the new column looks like it belongs, exactly the way a person
# would do it by hand in Excel (type the header, then copy the format of
# the neighbouring cell down).

$wb = $excel.ActiveWorkbook

# ---- TestSteps: new column F = "Run" ---------------------------------
$wsSteps = $wb.Worksheets.Item("TestSteps")

$wsSteps.Range("F1").Value = "Run"
$wsSteps.Range("E1").Copy()
$wsSteps.Range("F1").PasteSpecial(-4122)

$wsSteps.Range("F2:F4").Value = "Y"
$wsSteps.Range("E2:E4").Copy()
$wsSteps.Range("F2:F4").PasteSpecial(-4122)

$wsSteps.Range("F1").Select()

# ---- PageModules: new column I = "Run" --------------------------------
$wsMod = $wb.Worksheets.Item("PageModules")

$wsMod.Range("I1").Value = "Run"
$wsSteps.Range("F1").Copy()
$wsMod.Range("I1").PasteSpecial(-4122)

$wsMod.Range("I2:I5").Value = "Y"
$wsSteps.Range("F2:F4").Copy()
$wsMod.Range("I2:I4").PasteSpecial(-4122)
$wsSteps.Range("F2").Copy()
$wsMod.Range("I5").PasteSpecial(-4122)

$wsMod.Range("I8").Select()

# ---- TestCases: no data change, just leave the cursor on the existing
#      "Run" column (D) like the author did while reviewing it.
$wsCases = $wb.Worksheets.Item("TestCases")
$wsCases.Columns.Item("D").Select()
$wsCases.Range("D1").Activate()
